# Generate Report for Archive
# The localization status for "f6bf9a95-726d-4205-b211-4d668a1cd032.md"
# moved from "Ready for handoff" to "In Translation". Update this on the
# Overview sheet (zh-cn and de-de summary columns) as well as on the
# per-locale "zh-cn" and "de-de" sheets (Status column), for that file's row.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E8").Value = "In Translation"
$overview.Range("F8").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C8").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C8").Value = "In Translation"
